# Apply merged annotation updates to column E (label column) on Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Updated/merged labels in column E (row => new value)
$updates = @{
    5   = 1
    6   = 1
    9   = 1
    10  = 1
    11  = 1
    12  = 1
    13  = 1
    15  = 1
    23  = 1
    24  = 1
    28  = 1
    29  = 1
    33  = 1
    34  = 1
    36  = 1
    37  = 1
    38  = 1
    40  = 1
    41  = 1
    42  = 1
    45  = 1
    48  = 1
    50  = 1
    51  = 1
    55  = 1
    58  = 1
    61  = 0
    66  = 1
    78  = 1
    84  = 1
    91  = 1
    109 = 1
    116 = 1
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 5).Value = $updates[$row]
}

# Restore the view: scroll down to row 103 and leave the active selection on E116,
# matching where the editor was last working in the sheet.
$win = $excel.ActiveWindow
$ws.Range("E116").Select()
$win.ScrollRow = 103
$win.ScrollColumn = 1
